$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (index 1) from H1 to the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J (rows 2-11)
$iValues = @(10, 8, 9, 7, 8, 9, 7, 7, 6, 8)
$jValues = @(10, 9, 9, 8, 9, 9, 7, 8, 7, 9)

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jValues[$r - 2]
}
